{"js": "// Replace the date string and every \"A\u00f7B=C, D\" answer cell in the\n// practice table with the updated values from the commit.\n// Each old value is a unique substring in the document, so a\n// search-and-replace keyed on the exact old text is safe and keeps the\n// original run formatting (font/size) untouched, since only the\n// w:t text of the matched range is rewritten.\nconst replacements = [\n  [\"2025-08-12 Tuesday\", \"2025-08-13 Wednesday\"],\n  [\"918\u00f74=229, 2\", \"100\u00f73=33, 1\"],\n  [\"681\u00f72=340, 1\", \"588\u00f72=294, 0\"],\n  [\"148\u00f75=29, 3\", \"171\u00f73=57, 0\"],\n  [\"682\u00f73=227, 1\", \"118\u00f79=13, 1\"],\n  [\"430\u00f76=71, 4\", \"107\u00f77=15, 2\"],\n  [\"314\u00f77=44, 6\", \"781\u00f75=156, 1\"],\n  [\"335\u00f77=47, 6\", \"219\u00f75=43, 4\"],\n  [\"804\u00f72=402, 0\", \"232\u00f79=25, 7\"],\n  [\"935\u00f78=116, 7\", \"567\u00f73=189, 0\"],\n  [\"309\u00f77=44, 1\", \"308\u00f74=77, 0\"],\n  [\"328\u00f79=36, 4\", \"779\u00f79=86, 5\"],\n  [\"532\u00f73=177, 1\", \"653\u00f77=93, 2\"],\n  [\"364\u00f74=91, 0\", \"446\u00f72=223, 0\"],\n  [\"688\u00f74=172, 0\", \"731\u00f76=121, 5\"],\n  [\"436\u00f79=48, 4\", \"614\u00f79=68, 2\"],\n  [\"134\u00f75=26, 4\", \"586\u00f78=73, 2\"],\n  [\"595\u00f76=99, 1\", \"174\u00f78=21, 6\"],\n  [\"379\u00f74=94, 3\", \"661\u00f75=132, 1\"],\n  [\"149\u00f79=16, 5\", \"759\u00f74=189, 3\"],\n  [\"482\u00f76=80, 2\", \"507\u00f78=63, 3\"],\n  [\"549\u00f75=109, 4\", \"308\u00f72=154, 0\"],\n  [\"158\u00f74=39, 2\", \"572\u00f77=81, 5\"],\n  [\"333\u00f79=37, 0\", \"759\u00f76=126, 3\"],\n  [\"451\u00f77=64, 3\", \"605\u00f79=67, 2\"],\n  [\"380\u00f76=63, 2\", \"120\u00f77=17, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date string and every \"A\u00f7B=C, D\" answer cell in the\n# practice table with the updated values from the commit.\n# Each old value is a unique substring of the document, so a simple\n# Find/Replace (one match each, ReplaceOne) is safe and keeps the\n# original run formatting (font/size) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2025-08-12 Tuesday', '2025-08-13 Wednesday'),\n    @('918\u00f74=229, 2', '100\u00f73=33, 1'),\n    @('681\u00f72=340, 1', '588\u00f72=294, 0'),\n    @('148\u00f75=29, 3', '171\u00f73=57, 0'),\n    @('682\u00f73=227, 1', '118\u00f79=13, 1'),\n    @('430\u00f76=71, 4', '107\u00f77=15, 2'),\n    @('314\u00f77=44, 6', '781\u00f75=156, 1'),\n    @('335\u00f77=47, 6', '219\u00f75=43, 4'),\n    @('804\u00f72=402, 0', '232\u00f79=25, 7'),\n    @('935\u00f78=116, 7', '567\u00f73=189, 0'),\n    @('309\u00f77=44, 1', '308\u00f74=77, 0'),\n    @('328\u00f79=36, 4', '779\u00f79=86, 5'),\n    @('532\u00f73=177, 1', '653\u00f77=93, 2'),\n    @('364\u00f74=91, 0', '446\u00f72=223, 0'),\n    @('688\u00f74=172, 0', '731\u00f76=121, 5'),\n    @('436\u00f79=48, 4', '614\u00f79=68, 2'),\n    @('134\u00f75=26, 4', '586\u00f78=73, 2'),\n    @('595\u00f76=99, 1', '174\u00f78=21, 6'),\n    @('379\u00f74=94, 3', '661\u00f75=132, 1'),\n    @('149\u00f79=16, 5', '759\u00f74=189, 3'),\n    @('482\u00f76=80, 2', '507\u00f78=63, 3'),\n    @('549\u00f75=109, 4', '308\u00f72=154, 0'),\n    @('158\u00f74=39, 2', '572\u00f77=81, 5'),\n    @('333\u00f79=37, 0', '759\u00f76=126, 3'),\n    @('451\u00f77=64, 3', '605\u00f79=67, 2'),\n    @('380\u00f76=63, 2', '120\u00f77=17, 1')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue = 1, wdReplaceOne = 1\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 1)\n}\n"}
